# Atualização de bases das ligas, do dia: 20-06-2024 às 20:11
# Two pairs of rows had their match-data (everything except the leading
# sequence number in column A) swapped between each other:
#   row 104 <-> row 105
#   row 112 <-> row 113

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData {
    param($Row1, $Row2)

    # Columns B..AD hold the data that needs to be exchanged between the
    # two rows; column A (the running id) must stay untouched.
    $rng1 = $ws.Range("B$($Row1):AD$($Row1)")
    $rng2 = $ws.Range("B$($Row2):AD$($Row2)")

    $vals1 = $rng1.Value()
    $vals2 = $rng2.Value()

    $rng1.Value = $vals2
    $rng2.Value = $vals1
}

Swap-RowData 104 105
Swap-RowData 112 113
